$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 96; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq 71) {
        $cell.Value2 = 271
    }
}
